$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Capture the two existing comments (Email Recipients header row) before we
#        touch anything, so we can re-create them after the row shift below. ---
$commentA6 = $ws.Range("A6").Comment.Text()
$commentB6 = $ws.Range("B6").Comment.Text()

# --- 2. Re-order the config rows 2-4: "Spreadsheet URL" and "Sheet Name" move up,
#        "Email Subject" moves down to row 4. ---
$emailSubjectName  = $ws.Range("A2").Value()
$emailSubjectValue = $ws.Range("B2").Value()
$emailSubjectDesc  = $ws.Range("C2").Value()

$ws.Range("A2").Value = $ws.Range("A3").Value()
$ws.Range("B2").Value = $ws.Range("B3").Value()
$ws.Range("C2").Value = $ws.Range("C3").Value()

$ws.Range("A3").Value = $ws.Range("A4").Value()
$ws.Range("B3").Value = $ws.Range("B4").Value()
$ws.Range("C3").Value = $ws.Range("C4").Value()

$ws.Range("A4").Value = $emailSubjectName
$ws.Range("B4").Value = $emailSubjectValue
$ws.Range("C4").Value = $emailSubjectDesc

# --- 3. Remove the old comments - they will be re-added after the insert below
#        shifts their row down by one. ---
$ws.Range("A6").Comment.Delete()
$ws.Range("B6").Comment.Delete()

# --- 4. Insert a new blank row at row 5 for the "Email Text" entry. This pushes
#        the "Email Recipients" block (and its hyperlinks) down by one row. ---
$ws.Rows.Item(5).Insert()

# --- 5. Populate the new "Email Text" config row. ---
$ws.Range("A5").Value = "Email Text"
$ws.Range("B5").Value = "Mohon dapat disiapkan laptop dan email untuk Talents yang akan boarding dengan detail sebagai berikut:"
$ws.Range("C5").Value = "Text to be displayed in the body of the email"

# --- 6. Re-create the two comments on the "Email Recipients" header, now at row 7. ---
$ws.Range("A7").AddComment($commentA6)
$ws.Range("B7").AddComment($commentB6)

# --- 7. Update the active selection to match the edited workbook. ---
$ws.Range("A5").Select()
